$p = $ppt.ActivePresentation

$oldText = "https://www.udemy.com/course/angular-2-and-nodejs-the-practical-guide/learn/lecture/10419170#overview"
$newText = "https://github.com/peterhchen/900_MEAN_Proj"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    foreach ($shp in $s.Shapes) {
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldText) {
                $tr.Text = $newText
                $shp.TextFrame.AutoSize = 2
            }
        }
    }
}
